# language changing improved, currency class added and remembering

$wb = $excel.ActiveWorkbook

# --- Sheet1: refresh the latest rate snapshot values ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A1").Value = "29031,44"
$ws1.Range("A2").Value = "1831,15"
# force text so the comma-decimal reading (6 fractional digits) isn't
# mis-parsed as a grouped integer
$ws1.Range("A3").NumberFormat = "@"
$ws1.Range("A3").Value = "0,616699"

# --- data: add currency class headers + remembered sheet/cell refs ---
$wsData = $wb.Worksheets.Item("data")
$wsData.Range("A1").Value = "BITCOIN"
$wsData.Range("B1").Value = "ETHEREUM"
$wsData.Range("C1").Value = "XRP"

$wsData.Range("A2").Value = "Sheet1"
$wsData.Range("B2").Value = "Sheet1"
$wsData.Range("C2").Value = "Sheet1"

$wsData.Range("A3").Value = "A1"
$wsData.Range("B3").Value = "A2"
$wsData.Range("C3").Value = "A3"
